$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from the German default "Tabelle1" to "Sheet1"
$ws.Name = "Sheet1"

# The carrying_capacity_2040 and carrying_capacity_2050 columns (E and F)
# were dropped from the published data set; clear their header + values
# for rows 1-17 while leaving the per-cell number formatting intact.
$ws.Range("E1:F17").ClearContents()

# Update the remembered selection/scroll position to reflect where the
# editor was last working.
$ws.Range("D37").Select()
